$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.603177
$ws.Cells.Item(2, 8).Value = 19.809531
$ws.Cells.Item(2, 9).Value = 0.5135477412645301
$ws.Cells.Item(2, 10).Value = 0.5135477412645302
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1311436666666667
$ws.Cells.Item(2, 14).Value = 0.393431
$ws.Cells.Item(2, 15).Value = 0.02663441993971509
$ws.Cells.Item(2, 16).Value = 0.02663441993971509
$ws.Cells.Item(2, 17).Value = 0.8659648434289999
$ws.Cells.Item(2, 18).Value = 7.793683590861
$ws.Cells.Item(2, 19).Value = 0.01367804619993165
$ws.Cells.Item(2, 20).Value = 0.01367804619993165
$ws.Cells.Item(3, 7).Value = 6.603177
$ws.Cells.Item(3, 8).Value = 19.809531
$ws.Cells.Item(3, 9).Value = 0.5135477412645301
$ws.Cells.Item(3, 10).Value = 0.5135477412645302
$ws.Cells.Item(3, 15).Value = 0.06149297381279183
$ws.Cells.Item(3, 16).Value = 0.06149297381279183
$ws.Cells.Item(3, 17).Value = 1.999320937355
$ws.Cells.Item(3, 18).Value = 17.993888436195
$ws.Cells.Item(3, 19).Value = 0.03157957780519814
$ws.Cells.Item(3, 20).Value = 0.03157957780519815
$ws.Cells.Item(4, 7).Value = 6.603177
$ws.Cells.Item(4, 8).Value = 19.809531
$ws.Cells.Item(4, 9).Value = 0.5135477412645301
$ws.Cells.Item(4, 10).Value = 0.5135477412645302
$ws.Cells.Item(4, 13).Value = 4.009307333333333
$ws.Cells.Item(4, 14).Value = 12.027922
$ws.Cells.Item(4, 15).Value = 0.8142640654908683
$ws.Cells.Item(4, 16).Value = 0.8142640654908684
$ws.Cells.Item(4, 17).Value = 26.474165969398
$ws.Cells.Item(4, 18).Value = 238.267493724582
$ws.Cells.Item(4, 19).Value = 0.4181634716257088
$ws.Cells.Item(4, 20).Value = 0.4181634716257089
$ws.Cells.Item(5, 7).Value = 6.603177
$ws.Cells.Item(5, 8).Value = 19.809531
$ws.Cells.Item(5, 9).Value = 0.5135477412645301
$ws.Cells.Item(5, 10).Value = 0.5135477412645302
$ws.Cells.Item(5, 13).Value = 0.480609
$ws.Cells.Item(5, 14).Value = 1.441827
$ws.Cells.Item(5, 15).Value = 0.09760854075662465
$ws.Cells.Item(5, 16).Value = 0.09760854075662465
$ws.Cells.Item(5, 17).Value = 3.173546294793
$ws.Cells.Item(5, 18).Value = 28.561916653137
$ws.Cells.Item(5, 19).Value = 0.05012664563369141
$ws.Cells.Item(5, 20).Value = 0.05012664563369142
$ws.Cells.Item(6, 9).Value = 0.02944398858046029
$ws.Cells.Item(6, 10).Value = 0.0294439885804603
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1311436666666667
$ws.Cells.Item(6, 14).Value = 0.393431
$ws.Cells.Item(6, 15).Value = 0.02663441993971509
$ws.Cells.Item(6, 16).Value = 0.02663441993971509
$ws.Cells.Item(6, 17).Value = 0.04964963704877778
$ws.Cells.Item(6, 18).Value = 0.446846733439
$ws.Cells.Item(6, 19).Value = 0.000784223556552155
$ws.Cells.Item(6, 20).Value = 0.0007842235565521551
$ws.Cells.Item(7, 9).Value = 0.02944398858046029
$ws.Cells.Item(7, 10).Value = 0.0294439885804603
$ws.Cells.Item(7, 15).Value = 0.06149297381279183
$ws.Cells.Item(7, 16).Value = 0.06149297381279183
$ws.Cells.Item(7, 19).Value = 0.001810598418722386
$ws.Cells.Item(7, 20).Value = 0.001810598418722387
$ws.Cells.Item(8, 9).Value = 0.02944398858046029
$ws.Cells.Item(8, 10).Value = 0.0294439885804603
$ws.Cells.Item(8, 13).Value = 4.009307333333333
$ws.Cells.Item(8, 14).Value = 12.027922
$ws.Cells.Item(8, 15).Value = 0.8142640654908683
$ws.Cells.Item(8, 16).Value = 0.8142640654908684
$ws.Cells.Item(8, 17).Value = 1.517882326890889
$ws.Cells.Item(8, 18).Value = 13.660940942018
$ws.Cells.Item(8, 19).Value = 0.0239751818457923
$ws.Cells.Item(8, 20).Value = 0.0239751818457923
$ws.Cells.Item(9, 9).Value = 0.02944398858046029
$ws.Cells.Item(9, 10).Value = 0.0294439885804603
$ws.Cells.Item(9, 13).Value = 0.480609
$ws.Cells.Item(9, 14).Value = 1.441827
$ws.Cells.Item(9, 15).Value = 0.09760854075662465
$ws.Cells.Item(9, 16).Value = 0.09760854075662465
$ws.Cells.Item(9, 17).Value = 0.181953601107
$ws.Cells.Item(9, 18).Value = 1.637582409963
$ws.Cells.Item(9, 19).Value = 0.002873984759393449
$ws.Cells.Item(9, 20).Value = 0.002873984759393449
$ws.Cells.Item(10, 7).Value = 3.441487333333333
$ws.Cells.Item(10, 8).Value = 10.324462
$ws.Cells.Item(10, 9).Value = 0.2676541983690312
$ws.Cells.Item(10, 10).Value = 0.2676541983690313
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1311436666666667
$ws.Cells.Item(10, 14).Value = 0.393431
$ws.Cells.Item(10, 15).Value = 0.02663441993971509
$ws.Cells.Item(10, 16).Value = 0.02663441993971509
$ws.Cells.Item(10, 17).Value = 0.4513292676802222
$ws.Cells.Item(10, 18).Value = 4.061963409122
$ws.Cells.Item(10, 19).Value = 0.007128814317988583
$ws.Cells.Item(10, 20).Value = 0.007128814317988584
$ws.Cells.Item(11, 7).Value = 3.441487333333333
$ws.Cells.Item(11, 8).Value = 10.324462
$ws.Cells.Item(11, 9).Value = 0.2676541983690312
$ws.Cells.Item(11, 10).Value = 0.2676541983690313
$ws.Cells.Item(11, 15).Value = 0.06149297381279183
$ws.Cells.Item(11, 16).Value = 0.06149297381279183
$ws.Cells.Item(11, 17).Value = 1.042019270598889
$ws.Cells.Item(11, 18).Value = 9.37817343539
$ws.Cells.Item(11, 19).Value = 0.01645885261119063
$ws.Cells.Item(11, 20).Value = 0.01645885261119063
$ws.Cells.Item(12, 7).Value = 3.441487333333333
$ws.Cells.Item(12, 8).Value = 10.324462
$ws.Cells.Item(12, 9).Value = 0.2676541983690312
$ws.Cells.Item(12, 10).Value = 0.2676541983690313
$ws.Cells.Item(12, 13).Value = 4.009307333333333
$ws.Cells.Item(12, 14).Value = 12.027922
$ws.Cells.Item(12, 15).Value = 0.8142640654908683
$ws.Cells.Item(12, 16).Value = 0.8142640654908684
$ws.Cells.Item(12, 17).Value = 13.79798040310711
$ws.Cells.Item(12, 18).Value = 124.181823627964
$ws.Cells.Item(12, 19).Value = 0.2179411957096667
$ws.Cells.Item(12, 20).Value = 0.2179411957096667
$ws.Cells.Item(13, 7).Value = 3.441487333333333
$ws.Cells.Item(13, 8).Value = 10.324462
$ws.Cells.Item(13, 9).Value = 0.2676541983690312
$ws.Cells.Item(13, 10).Value = 0.2676541983690313
$ws.Cells.Item(13, 13).Value = 0.480609
$ws.Cells.Item(13, 14).Value = 1.441827
$ws.Cells.Item(13, 15).Value = 0.09760854075662465
$ws.Cells.Item(13, 16).Value = 0.09760854075662465
$ws.Cells.Item(13, 17).Value = 1.654009785786
$ws.Cells.Item(13, 18).Value = 14.886088072074
$ws.Cells.Item(13, 19).Value = 0.02612533573018528
$ws.Cells.Item(13, 20).Value = 0.02612533573018529
$ws.Cells.Item(14, 7).Value = 2.434707333333333
$ws.Cells.Item(14, 8).Value = 7.304122
$ws.Cells.Item(14, 9).Value = 0.1893540717859783
$ws.Cells.Item(14, 10).Value = 0.1893540717859783
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.1311436666666667
$ws.Cells.Item(14, 14).Value = 0.393431
$ws.Cells.Item(14, 15).Value = 0.02663441993971509
$ws.Cells.Item(14, 16).Value = 0.02663441993971509
$ws.Cells.Item(14, 17).Value = 0.3192964469535555
$ws.Cells.Item(14, 18).Value = 2.873668022582
$ws.Cells.Item(14, 19).Value = 0.005043335865242702
$ws.Cells.Item(14, 20).Value = 0.005043335865242703
$ws.Cells.Item(15, 7).Value = 2.434707333333333
$ws.Cells.Item(15, 8).Value = 7.304122
$ws.Cells.Item(15, 9).Value = 0.1893540717859783
$ws.Cells.Item(15, 10).Value = 0.1893540717859783
$ws.Cells.Item(15, 15).Value = 0.06149297381279183
$ws.Cells.Item(15, 16).Value = 0.06149297381279183
$ws.Cells.Item(15, 17).Value = 0.7371847442322222
$ws.Cells.Item(15, 18).Value = 6.634662698090001
$ws.Cells.Item(15, 19).Value = 0.01164394497768067
$ws.Cells.Item(15, 20).Value = 0.01164394497768067
$ws.Cells.Item(16, 7).Value = 2.434707333333333
$ws.Cells.Item(16, 8).Value = 7.304122
$ws.Cells.Item(16, 9).Value = 0.1893540717859783
$ws.Cells.Item(16, 10).Value = 0.1893540717859783
$ws.Cells.Item(16, 13).Value = 4.009307333333333
$ws.Cells.Item(16, 14).Value = 12.027922
$ws.Cells.Item(16, 15).Value = 0.8142640654908683
$ws.Cells.Item(16, 16).Value = 0.8142640654908684
$ws.Cells.Item(16, 17).Value = 9.761489966053777
$ws.Cells.Item(16, 18).Value = 87.853409694484
$ws.Cells.Item(16, 19).Value = 0.1541842163097004
$ws.Cells.Item(16, 20).Value = 0.1541842163097004
$ws.Cells.Item(17, 7).Value = 2.434707333333333
$ws.Cells.Item(17, 8).Value = 7.304122
$ws.Cells.Item(17, 9).Value = 0.1893540717859783
$ws.Cells.Item(17, 10).Value = 0.1893540717859783
$ws.Cells.Item(17, 13).Value = 0.480609
$ws.Cells.Item(17, 14).Value = 1.441827
$ws.Cells.Item(17, 15).Value = 0.09760854075662465
$ws.Cells.Item(17, 16).Value = 0.09760854075662465
$ws.Cells.Item(17, 17).Value = 1.170142256766
$ws.Cells.Item(17, 18).Value = 10.531280310894
$ws.Cells.Item(17, 19).Value = 0.01848257463335449
$ws.Cells.Item(17, 20).Value = 0.0184825746333545
